# Ajuste no cabeçalho da planilha
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing header cells to the new standardized (upper-case / underscored) labels ---
$ws.Range("B1").Value = "NOME_FANTASIA"
$ws.Range("C1").Value = "ENDERECO"
$ws.Range("D1").Value = "NUMERO"
$ws.Range("I1").Value = "TIPO_TITULO"
$ws.Range("K1").Value = "PLANO"
$ws.Range("L1").Value = "CONTRATO"
$ws.Range("M1").Value = "DOCUMENTO"
$ws.Range("N1").Value = "PRODUTO"
$ws.Range("Q1").Value = "CELULAR_WHATSAPP"
$ws.Range("R1").Value = "FONE_COMERCIAL"
$ws.Range("S1").Value = "FONE_RESIDENCIAL"
$ws.Range("T1").Value = "FONE_OUTROS"

# --- Type the two new trailing header labels (plain cells first) ---
$ws.Range("U1").Value = "EMAIL_1"
$ws.Range("V1").Value = "EMAIL_2"

# --- Underline the whole header row (kept bold) ---
$ws.Range("A1:V1").Font.Underline = 2   # xlUnderlineStyleSingle

# --- Bring the two new cells up to the same look as the rest of the header
#     (bold/underline font + fill + border), by copying the last header cell's format ---
$ws.Range("T1").Copy() | Out-Null
$ws.Range("U1:V1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Re-fit the header row now that the label text has changed ---
$ws.Range("A1:T1").EntireColumn.AutoFit() | Out-Null

# --- Update the worksheet selection to the full header row (matches the authored state) ---
$ws.Range("A1:XFD1").Select() | Out-Null
